# 200. Number of Islands
# Append a new tracker row (row 99) to Sheet1, mirroring the formatting of
# the previous entry row (row 98), then fill in the new problem's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 98 (A:E) down to row 99 first so the new row inherits the same
# cell formatting/styles used by the other entry rows (bold-ish link style
# in column B, date style in column E, etc.).
$ws.Range("A98:E98").Copy($ws.Range("A99:E99"))

# Overwrite the copied values with the new entry's data.
$ws.Range("A99").Value = 200
$ws.Range("B99").Value = "Number of Islands"
$ws.Range("C99").Value = "Medium"
$ws.Range("D99").Value = "Graph , BFS , queue ,counting"
$ws.Range("E99").Value = 45806

# Update the view: selection moves to B97 and the view is scrolled down so
# row 80 is at the top.
$null = $ws.Range("B97").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 80
$win.ScrollColumn = 1
